$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the 2019 row (row 2) dates - they were mistakenly duplicated from 2024
$ws.Range("B2").Value = 43764
$ws.Range("C2").Value = 43800
$ws.Range("D2").Value = 43786

# Update the active cell selection as recorded in the saved file
$ws.Range("I9").Select()
